$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@(2, '25.999.62', '  +0.17%  ')
    ,@(3, '1.641.55', '  -0.53%  ')
    ,@(4, '1.001', '  -0.97%  ')
    ,@(5, '215.14', '  -0.54%  ')
    ,@(6, '0.5062', '  -0.98%  ')
    ,@(7, $null, '  -0.68%  ')
    ,@(8, '0.2579', '  -0.25%  ')
    ,@(9, '0.06358', '  -1.17%  ')
    ,@(10, '19.85', '  +0.77%  ')
    ,@(11, '0.07758', '  -0.44%  ')
    ,@(12, '4.281', '  -1.14%  ')
    ,@(13, '1.633.58', '  -0.98%  ')
    ,@(14, '0.5484', '  +0.12%  ')
    ,@(15, '0.0₅7749', '  -1.93%  ')
    ,@(16, '64.39', '  -0.75%  ')
    ,@(17, '26.022.28', '  -0.08%  ')
    ,@(18, '1.002', '  -0.75%  ')
    ,@(19, '196.52', '  -1.06%  ')
    ,@(20, '4.441', '  -0.82%  ')
    ,@(21, '9.949', '  -0.90%  ')
    ,@(22, '6.105', '  +0.32%  ')
    ,@(23, $null, '  -0.58%  ')
    ,@(24, $null, '  +1.66%  ')
    ,@(25, '143.84', '  +2.37%  ')
    ,@(26, '0.1236', '  +7.23%  ')
    ,@(27, '6.877', $null)
    ,@(28, '15.65', '  -0.85%  ')
    ,@(29, $null, '  -0.53%  ')
    ,@(30, '0.04875', '  -3.05%  ')
    ,@(31, '3.273', '  -0.58%  ')
    ,@(32, '3.223', '  +0.38%  ')
    ,@(33, '1.545', '  -0.24%  ')
    ,@(34, '2.376', '  +0.26%  ')
    ,@(35, '0.9153', '  +2.07%  ')
    ,@(36, $null, '  -1.19%  ')
    ,@(37, '0.5545', $null)
    ,@(38, '1.090.33', '  -4.18%  ')
    ,@(39, '0.01572', $null)
    ,@(40, '1.002', '  -0.81%  ')
    ,@(41, '2.525', '  -1.81%  ')
    ,@(42, '5.610', '  -1.13%  ')
    ,@(43, '0.8064', '  -1.57%  ')
    ,@(44, '99.19', '  -0.94%  ')
    ,@(45, $null, '  -4.39%  ')
    ,@(46, '1.778.94', '  -0.43%  ')
    ,@(47, '0.4537', '  -0.35%  ')
    ,@(48, '55.46', '  +0.17%  ')
    ,@(49, '1.003', '  -0.52%  ')
    ,@(50, '0.05216', '  +2.23%  ')
    ,@(51, '7.537', '  +1.92%  ')
)

foreach ($u in $updates) {
    $row = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]
    if ($dVal -ne $null) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.Value = "'" + $dVal
        $cell.Style = "Normal"
    }
    if ($eVal -ne $null) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}
